$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 173 (pushes old rows 173:185 down to 174:186),
# inheriting formatting (date style on column D) from the row below, like
# Excel's native Insert behaviour.
$ws.Rows(173).Insert()

# Populate the newly inserted row 173 with the new weekly record.
$ws.Range("A173").Value = 10
$ws.Range("B173").Value = "Vega Modelo de Temuco"
$ws.Range("C173").Value = "La Araucanía"
$ws.Range("D173").Value = 45013
$ws.Range("E173").Value = 9
$ws.Range("F173").Value = 100114002
$ws.Range("G173").Value = "Camote"
$ws.Range("H173").Value = "Sin especificar"
$ws.Range("I173").Value = "Primera"
$ws.Range("J173").Value = 40
$ws.Range("K173").Value = 20000
$ws.Range("L173").Value = 20000
$ws.Range("M173").Value = 20000
$ws.Range("N173").Value = "$/malla 20 kilos"
$ws.Range("O173").Value = "Perú"
$ws.Range("P173").Value = 1000
$ws.Range("Q173").Value = 20
$ws.Range("R173").Value = "Hortaliza"
